# Auto-generated cell updates for Pandaemonium_Profits workbook
# Applies numeric corrections to columns H-N across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 166.45
$ws.Range("I15").Value = 166.45
$ws.Range("K15").Value = 499.35
$ws.Range("M15").Value = -330.35
# Row 97
$ws.Range("H97").Value = 0.0
$ws.Range("J97").Value = 0.0
$ws.Range("L97").Value = 0.0
$ws.Range("N97").ClearContents()
# Row 112
$ws.Range("H112").Value = 1817.05
$ws.Range("J112").Value = 1974.5
$ws.Range("L112").Value = 5923.5
$ws.Range("N112").Value = -8139.5
# Row 135
$ws.Range("H135").Value = 45457460.0
$ws.Range("I135").Value = 17859148.0
$ws.Range("J135").Value = 200008000.0
$ws.Range("K135").Value = 160732332.0
$ws.Range("L135").Value = 1800072000.0
$ws.Range("M135").Value = -160729797.0
$ws.Range("N135").Value = -1800077070.0
# Row 137
$ws.Range("H137").Value = 2744.611
$ws.Range("I137").Value = 1495.5927
$ws.Range("J137").Value = 3993.6296
$ws.Range("K137").Value = 4486.7781
$ws.Range("L137").Value = 11980.8888
$ws.Range("M137").Value = -1936.7781
$ws.Range("N137").Value = -17080.8888
# Row 138
$ws.Range("H138").Value = 1572826.1
$ws.Range("J138").Value = 1753656.5
$ws.Range("L138").Value = 5260969.5
$ws.Range("N138").Value = -5271249.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23059.572
$ws.Range("I32").Value = 17260.133
$ws.Range("J32").Value = 26281.482
$ws.Range("K32").Value = 17260.133
$ws.Range("L32").Value = 26281.482
$ws.Range("M32").Value = -16973.133
$ws.Range("N32").Value = -26855.482
# Row 122
$ws.Range("H122").Value = 7814613.5
$ws.Range("I122").Value = 1829.0
$ws.Range("K122").Value = 5487.0
$ws.Range("M122").Value = -3037.0
# Row 132
$ws.Range("H132").Value = 4446.9814
$ws.Range("I132").Value = 1758.7241
$ws.Range("J132").Value = 7565.36
$ws.Range("K132").Value = 5276.1723
$ws.Range("L132").Value = 22696.08
$ws.Range("M132").Value = -2746.1723
$ws.Range("N132").Value = -27756.08

$ws = $wb.Worksheets.Item("BSM")
# Row 44
$ws.Range("H44").Value = 15000.0
$ws.Range("J44").Value = 15000.0
$ws.Range("L44").Value = 15000.0
$ws.Range("N44").Value = -15994.0
# Row 132
$ws.Range("H132").Value = 61561.25
$ws.Range("J132").Value = 61561.25
$ws.Range("L132").Value = 61561.25
$ws.Range("N132").Value = -71681.25
# Row 134
$ws.Range("H134").Value = 21950.71
$ws.Range("I134").Value = 2386.725
$ws.Range("J134").Value = 87164.0
$ws.Range("K134").Value = 7160.174999999999
$ws.Range("L134").Value = 261492.0
$ws.Range("M134").Value = -4625.174999999999
$ws.Range("N134").Value = -266562.0

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 37455.0
$ws.Range("I23").Value = 30000.0
$ws.Range("J23").Value = 39940.0
$ws.Range("K23").Value = 30000.0
$ws.Range("L23").Value = 39940.0
$ws.Range("M23").Value = -29760.0
$ws.Range("N23").Value = -40420.0
# Row 27
$ws.Range("H27").Value = 37455.0
$ws.Range("I27").Value = 30000.0
$ws.Range("J27").Value = 39940.0
$ws.Range("K27").Value = 30000.0
$ws.Range("L27").Value = 39940.0
$ws.Range("M27").Value = -29808.0
$ws.Range("N27").Value = -40324.0
# Row 31
$ws.Range("H31").Value = 4111.4
$ws.Range("I31").Value = 1570.2632
$ws.Range("J31").Value = 5968.385
$ws.Range("K31").Value = 1570.2632
$ws.Range("L31").Value = 5968.385
$ws.Range("M31").Value = -1275.2632
$ws.Range("N31").Value = -6558.385
# Row 34
$ws.Range("H34").Value = 4111.4
$ws.Range("I34").Value = 1570.2632
$ws.Range("J34").Value = 5968.385
$ws.Range("K34").Value = 1570.2632
$ws.Range("L34").Value = 5968.385
$ws.Range("M34").Value = -1368.2632
$ws.Range("N34").Value = -6372.385
# Row 135
$ws.Range("H135").Value = 58000.0
$ws.Range("J135").Value = 58000.0
$ws.Range("L135").Value = 58000.0
$ws.Range("N135").Value = -68140.0

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 8574.875
$ws.Range("J39").Value = 8574.875
$ws.Range("L39").Value = 25724.625
$ws.Range("N39").Value = -26312.625
# Row 107
$ws.Range("H107").Value = 2778926.5
$ws.Range("J107").Value = 1324.0938
$ws.Range("L107").Value = 3972.2814
$ws.Range("N107").Value = -7812.2814
# Row 122
$ws.Range("H122").Value = 1455.0
$ws.Range("J122").Value = 1867.8
$ws.Range("L122").Value = 16810.2
$ws.Range("N122").Value = -21710.2

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 12000.0
$ws.Range("I12").Value = 0.0
$ws.Range("J12").Value = 12000.0
$ws.Range("K12").Value = 0.0
$ws.Range("L12").Value = 12000.0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -12280.0
# Row 102
$ws.Range("H102").Value = 6486.1875
$ws.Range("I102").Value = 6072.636
$ws.Range("K102").Value = 6072.636
$ws.Range("M102").Value = -4450.636
# Row 122
$ws.Range("H122").Value = 7977.7144
$ws.Range("I122").Value = 7243.222
$ws.Range("J122").Value = 9299.8
$ws.Range("K122").Value = 21729.666
$ws.Range("L122").Value = 27899.4
$ws.Range("M122").Value = -19279.666
$ws.Range("N122").Value = -32799.39999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5118.864
$ws.Range("I7").Value = 3159.4119
$ws.Range("J7").Value = 11781.0
$ws.Range("K7").Value = 3159.4119
$ws.Range("L7").Value = 11781.0
$ws.Range("M7").Value = -3047.4119
$ws.Range("N7").Value = -12005.0
# Row 38
$ws.Range("H38").Value = 19999.0
$ws.Range("J38").Value = 19999.0
$ws.Range("L38").Value = 19999.0
$ws.Range("N38").Value = -20819.0
# Row 40
$ws.Range("H40").Value = 5060.9
$ws.Range("I40").Value = 4938.0
$ws.Range("K40").Value = 4938.0
$ws.Range("M40").Value = -4802.0
# Row 61
$ws.Range("H61").Value = 919560.06
$ws.Range("I61").Value = 1378759.2
$ws.Range("J61").Value = 1161.875
$ws.Range("K61").Value = 1378759.2
$ws.Range("L61").Value = 1161.875
$ws.Range("M61").Value = -1378557.2
$ws.Range("N61").Value = -1565.875
# Row 74
$ws.Range("H74").Value = 36217.0
$ws.Range("J74").Value = 36217.0
$ws.Range("L74").Value = 36217.0
$ws.Range("N74").Value = -38213.0
# Row 77
$ws.Range("H77").Value = 36217.0
$ws.Range("J77").Value = 36217.0
$ws.Range("L77").Value = 108651.0
$ws.Range("N77").Value = -118635.0
# Row 113
$ws.Range("H113").Value = 919560.06
$ws.Range("I113").Value = 1378759.2
$ws.Range("J113").Value = 1161.875
$ws.Range("K113").Value = 1378759.2
$ws.Range("L113").Value = 1161.875
$ws.Range("M113").Value = -1376589.2
$ws.Range("N113").Value = -5501.875
# Row 122
$ws.Range("H122").Value = 7220.8096
$ws.Range("I122").Value = 5450.2856
$ws.Range("J122").Value = 10761.857
$ws.Range("K122").Value = 16350.8568
$ws.Range("L122").Value = 32285.571
$ws.Range("M122").Value = -13900.8568
$ws.Range("N122").Value = -37185.571
# Row 126
$ws.Range("H126").Value = 5118.864
$ws.Range("I126").Value = 3159.4119
$ws.Range("J126").Value = 11781.0
$ws.Range("K126").Value = 9478.235700000001
$ws.Range("L126").Value = 35343.0
$ws.Range("M126").Value = -7008.235700000001
$ws.Range("N126").Value = -40283.0

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1216.75
$ws.Range("I100").Value = 372.08334
$ws.Range("J100").Value = 3750.75
$ws.Range("K100").Value = 744.16668
$ws.Range("L100").Value = 7501.5
$ws.Range("M100").Value = -203.16668
$ws.Range("N100").Value = -8583.5
# Row 113
$ws.Range("H113").Value = 9728.0
$ws.Range("I113").Value = 1680.0
$ws.Range("J113").Value = 15476.571
$ws.Range("K113").Value = 5040.0
$ws.Range("L113").Value = 46429.713
$ws.Range("M113").Value = -2870.0
$ws.Range("N113").Value = -50769.713
# Row 126
$ws.Range("H126").Value = 1284.3684
$ws.Range("I126").Value = 1143.4286
$ws.Range("J126").Value = 1679.0
$ws.Range("K126").Value = 3430.2858
$ws.Range("L126").Value = 5037.0
$ws.Range("M126").Value = -960.2857999999997
$ws.Range("N126").Value = -9977.0
# Row 135
$ws.Range("H135").Value = 142885310.0
$ws.Range("J135").Value = 142885310.0
$ws.Range("L135").Value = 142885310.0
$ws.Range("N135").Value = -142895450.0
# Row 136
$ws.Range("H136").Value = 5632.5083
$ws.Range("I136").Value = 3212.5151
$ws.Range("J136").Value = 8704.038
$ws.Range("K136").Value = 9637.5453
$ws.Range("L136").Value = 26112.114
$ws.Range("M136").Value = -7087.5453
$ws.Range("N136").Value = -31212.114
# Row 137
$ws.Range("H137").Value = 60482.5
$ws.Range("J137").Value = 60482.5
$ws.Range("L137").Value = 60482.5
$ws.Range("N137").Value = -70682.5
# Row 139
$ws.Range("H139").Value = 68370.71
$ws.Range("J139").Value = 68370.71
$ws.Range("L139").Value = 68370.71
$ws.Range("N139").Value = -78650.71
